$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.080.34"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.96"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.51"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6356"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.86"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.04%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07346"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2934"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.79"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.829.66"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6640"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.07"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008676"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.057"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.080.67"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.076.15"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.45"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.44"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.484"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1363"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.92"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.504"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.090"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.203"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05329"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7372"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.650"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.301.59"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01788"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.744"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.317"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9029"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.62"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.975.13"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5134"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.15"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.731"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05815"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.94%  "
